# feat: add 2022-Q4 data
#
# Before:
#   Sheet1 "总计"    -> totals table (row2 = 2022-Q3 summary)
#   Sheet2 "2022-Q3" -> fund holdings detail for 2022-Q3 (13 data rows)
#
# After:
#   Sheet1 "总计"    -> totals table (row2 = 2022-Q4 summary, row3 = 2022-Q3 summary)
#   Sheet2 "2022-Q4" -> fund holdings detail for 2022-Q4 (2 data rows) [was "2022-Q3"]
#   Sheet3 "2022-Q3" -> fund holdings detail for 2022-Q3 (13 data rows, moved to a new sheet)

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item(1)
$q3Sheet = $wb.Worksheets.Item(2)
$styleSrc = $totalSheet.Range("A2")

# 1. Create a brand-new sheet right after the existing "2022-Q3" sheet. This
#    new sheet will keep the full 2022-Q3 fund-holdings detail (13 rows),
#    freeing up the original sheet2 to become the 2022-Q4 detail sheet.
$newQ3Sheet = $wb.Worksheets.Add($null, $q3Sheet)

# Copy the entire existing 2022-Q3 detail (values + styles) onto the new sheet
# before we overwrite the original sheet2 with the 2022-Q4 numbers.
$q3Sheet.Range("A1:H13").Copy($newQ3Sheet.Range("A1:H13"))

# Rename the ORIGINAL sheet first so the "2022-Q3" name is freed up before we
# try to assign it to the new sheet (Excel refuses duplicate sheet names).
$q3Sheet.Name = "2022-Q4"
$newQ3Sheet.Name = "2022-Q3"

# 2. Turn the original sheet2 (now named "2022-Q4") into the 2022-Q4 detail
#    sheet: drop the old 2022-Q3 rows 4-13 (only 2 funds were held in
#    2022-Q4).
$q3Sheet.Range("A4:H13").EntireRow.Delete()

# Header row: keep the bold/centered "index" style (style used by A2 on the
# totals sheet) instead of the old plain header style.
$styleSrc.Copy($q3Sheet.Range("B1"))
$q3Sheet.Range("B1").Value = "基金代码"
$styleSrc.Copy($q3Sheet.Range("C1"))
$q3Sheet.Range("C1").Value = "基金名称"
$styleSrc.Copy($q3Sheet.Range("D1"))
$q3Sheet.Range("D1").Value = "基金规模"
$styleSrc.Copy($q3Sheet.Range("E1"))
$q3Sheet.Range("E1").Value = "股票总仓位"
$styleSrc.Copy($q3Sheet.Range("F1"))
$q3Sheet.Range("F1").Value = "仓位占比"
$styleSrc.Copy($q3Sheet.Range("G1"))
$q3Sheet.Range("G1").Value = "持有市值(亿元)"
$styleSrc.Copy($q3Sheet.Range("H1"))
$q3Sheet.Range("H1").Value = "仓位排名"

# Row 2 - 南方宝丰混合A
$styleSrc.Copy($q3Sheet.Range("A2"))
$q3Sheet.Range("A2").Value = 0
$q3Sheet.Range("B2").NumberFormat = "@"
$q3Sheet.Range("B2").Value = "008513"
$q3Sheet.Range("C2").Value = "南方宝丰混合A"
$q3Sheet.Range("D2").NumberFormat = "@"
$q3Sheet.Range("D2").Value = "33.10"
$q3Sheet.Range("E2").NumberFormat = "@"
$q3Sheet.Range("E2").Value = "27.79"
$q3Sheet.Range("F2").NumberFormat = "@"
$q3Sheet.Range("F2").Value = "0.58"
$q3Sheet.Range("G2").NumberFormat = "@"
$q3Sheet.Range("G2").Value = "0.1920"
$q3Sheet.Range("H2").Value = 10

# Row 3 - 南方宝丰混合C
$styleSrc.Copy($q3Sheet.Range("A3"))
$q3Sheet.Range("A3").Value = 1
$q3Sheet.Range("B3").NumberFormat = "@"
$q3Sheet.Range("B3").Value = "008514"
$q3Sheet.Range("C3").Value = "南方宝丰混合C"
$q3Sheet.Range("D3").NumberFormat = "@"
$q3Sheet.Range("D3").Value = "3.34"
$q3Sheet.Range("E3").NumberFormat = "@"
$q3Sheet.Range("E3").Value = "27.79"
$q3Sheet.Range("F3").NumberFormat = "@"
$q3Sheet.Range("F3").Value = "0.58"
$q3Sheet.Range("G3").NumberFormat = "@"
$q3Sheet.Range("G3").Value = "0.0194"
$q3Sheet.Range("H3").Value = 10

# 3. Update the "总计" summary sheet: row 2 now reports 2022-Q4 figures, and a
#    new row 3 carries forward the old 2022-Q3 figures that used to live in
#    row 2.
$totalSheet.Range("A2:D2").Copy($totalSheet.Range("A3:D3"))
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q3"
$totalSheet.Range("C3").Value = 12
$totalSheet.Range("D3").Value = 0.78

$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.21

Write-Output "applied 2022-Q4 update"
